$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 1000
$ws.Range("B4").Value = 2000
$ws.Range("C4").Value = 3000

$ws.Range("A5").Value = 10000
$ws.Range("B5").Value = 20000
$ws.Range("C5").Value = 30000

$ws.Range("F10").Select()
